# Auto-generated Excel COM-interop edit script
# Applies value updates to Sheets ALC, ARM, BSM, CRP, GSM, LTW, WVR
# as described by the commit "chore: update Sheets via scheduled runner"

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 658497.75
$ws.Range("I98").Value = 745864.4
$ws.Range("J98").Value = 3248
$ws.Range("K98").Value = 745864.4
$ws.Range("L98").Value = 3248
$ws.Range("M98").Value = -744366.4
$ws.Range("N98").Value = -6244
# Row 122
$ws.Range("H122").Value = 658497.75
$ws.Range("I122").Value = 745864.4
$ws.Range("J122").Value = 3248
$ws.Range("K122").Value = 2237593.2
$ws.Range("L122").Value = 9744
$ws.Range("M122").Value = -2235143.2
$ws.Range("N122").Value = -14644
# Row 137
$ws.Range("H137").Value = 1688.2174
$ws.Range("I137").Value = 1757.2222
$ws.Range("K137").Value = 5271.6666
$ws.Range("M137").Value = -2721.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2032.4667
$ws.Range("I61").Value = 1297.1666
$ws.Range("J61").Value = 4973.6665
$ws.Range("K61").Value = 1297.1666
$ws.Range("L61").Value = 4973.6665
$ws.Range("M61").Value = -1085.1666
$ws.Range("N61").Value = -5397.6665
# Row 74
$ws.Range("H74").Value = 6613.091
$ws.Range("I74").Value = 1240
$ws.Range("K74").Value = 1240
$ws.Range("M74").Value = -366
# Row 77
$ws.Range("H77").Value = 6613.091
$ws.Range("I77").Value = 1240
$ws.Range("K77").Value = 6200
$ws.Range("M77").Value = -1832
# Row 125
$ws.Range("H125").Value = 34000
$ws.Range("J125").Value = 34000
$ws.Range("L125").Value = 34000
$ws.Range("N125").Value = -43840
# Row 136
$ws.Range("H136").Value = 2032.4667
$ws.Range("I136").Value = 1297.1666
$ws.Range("J136").Value = 4973.6665
$ws.Range("K136").Value = 3891.4998
$ws.Range("L136").Value = 14920.9995
$ws.Range("M136").Value = -1341.4998
$ws.Range("N136").Value = -20020.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 6816.7144
$ws.Range("I86").Value = 2065.8
$ws.Range("J86").Value = 18694
$ws.Range("K86").Value = 2065.8
$ws.Range("L86").Value = 18694
$ws.Range("M86").Value = -942.8000000000002
$ws.Range("N86").Value = -20940
# Row 89
$ws.Range("H89").Value = 6816.7144
$ws.Range("I89").Value = 2065.8
$ws.Range("J89").Value = 18694
$ws.Range("K89").Value = 10329
$ws.Range("L89").Value = 93470
$ws.Range("M89").Value = -4713
$ws.Range("N89").Value = -104702
# Row 134
$ws.Range("H134").Value = 3217.96
$ws.Range("I134").Value = 2347.7368
$ws.Range("J134").Value = 5973.6665
$ws.Range("K134").Value = 7043.2104
$ws.Range("L134").Value = 17920.9995
$ws.Range("M134").Value = -4508.2104
$ws.Range("N134").Value = -22990.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 37500
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25480
# Row 27
$ws.Range("H27").Value = 37500
$ws.Range("J27").Value = 25000
$ws.Range("L27").Value = 25000
$ws.Range("N27").Value = -25384
# Row 31
$ws.Range("H31").Value = 1682.0333
$ws.Range("I31").Value = 1090.4642
$ws.Range("K31").Value = 1090.4642
$ws.Range("M31").Value = -795.4641999999999
# Row 34
$ws.Range("H34").Value = 1682.0333
$ws.Range("I34").Value = 1090.4642
$ws.Range("K34").Value = 1090.4642
$ws.Range("M34").Value = -888.4641999999999
# Row 58
$ws.Range("H58").Value = 1789.12
$ws.Range("I58").Value = 976.35297
$ws.Range("J58").Value = 3516.25
$ws.Range("K58").Value = 976.35297
$ws.Range("L58").Value = 3516.25
$ws.Range("M58").Value = -773.35297
$ws.Range("N58").Value = -3922.25
# Row 132
$ws.Range("H132").Value = 3058.9614
$ws.Range("I132").Value = 1963.4117
$ws.Range("J132").Value = 5128.3335
$ws.Range("K132").Value = 5890.2351
$ws.Range("L132").Value = 15385.0005
$ws.Range("M132").Value = -3360.2351
$ws.Range("N132").Value = -20445.0005
# Row 134
$ws.Range("H134").Value = 3869.1738
$ws.Range("I134").Value = 2425.1875
$ws.Range("J134").Value = 7169.7144
$ws.Range("K134").Value = 7275.5625
$ws.Range("L134").Value = 21509.1432
$ws.Range("M134").Value = -4740.5625
$ws.Range("N134").Value = -26579.1432
# Row 136
$ws.Range("H136").Value = 1789.12
$ws.Range("I136").Value = 976.35297
$ws.Range("J136").Value = 3516.25
$ws.Range("K136").Value = 2929.05891
$ws.Range("L136").Value = 10548.75
$ws.Range("M136").Value = -379.0589100000002
$ws.Range("N136").Value = -15648.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 4
$ws.Range("H4").Value = 70004
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = $null
# Row 92
$ws.Range("H92").Value = 7061
$ws.Range("J92").Value = 7061
$ws.Range("L92").Value = 7061
$ws.Range("N92").Value = -10805
# Row 98
$ws.Range("H98").Value = 30000
$ws.Range("J98").Value = 30000
$ws.Range("L98").Value = 30000
$ws.Range("N98").Value = -35990
# Row 101
$ws.Range("H101").Value = 29000
$ws.Range("J101").Value = 29000
$ws.Range("L101").Value = 29000
$ws.Range("N101").Value = -35490
# Row 107
$ws.Range("H107").Value = 766.8823
$ws.Range("J107").Value = 410.42856
$ws.Range("L107").Value = 410.42856
$ws.Range("N107").Value = -4250.42856
# Row 126
$ws.Range("H126").Value = 2205.4878
$ws.Range("I126").Value = 1896.2142
$ws.Range("J126").Value = 2365.8518
$ws.Range("K126").Value = 5688.642599999999
$ws.Range("L126").Value = 7097.555399999999
$ws.Range("M126").Value = -3218.642599999999
$ws.Range("N126").Value = -12037.5554

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 95
$ws.Range("H95").Value = 59995
$ws.Range("J95").Value = 59995
$ws.Range("L95").Value = 59995
$ws.Range("N95").Value = -65487
# Row 97
$ws.Range("H97").Value = 14500
$ws.Range("J97").Value = 14500
$ws.Range("L97").Value = 14500
$ws.Range("N97").Value = -16482
# Row 132
$ws.Range("H132").Value = 4070.5334
$ws.Range("I132").Value = 2433.7646
$ws.Range("J132").Value = 6210.923
$ws.Range("K132").Value = 7301.293799999999
$ws.Range("L132").Value = 18632.769
$ws.Range("M132").Value = -4771.293799999999
$ws.Range("N132").Value = -23692.769
# Row 136
$ws.Range("H136").Value = 2655.3948
$ws.Range("I136").Value = 1413.25
$ws.Range("J136").Value = 3558.7727
$ws.Range("K136").Value = 4239.75
$ws.Range("L136").Value = 10676.3181
$ws.Range("M136").Value = -1689.75
$ws.Range("N136").Value = -15776.3181

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 15800
$ws.Range("J22").Value = 15800
$ws.Range("L22").Value = 15800
$ws.Range("N22").Value = -16386
# Row 92
$ws.Range("H92").Value = 49997.5
$ws.Range("J92").Value = 49997.5
$ws.Range("L92").Value = 49997.5
$ws.Range("N92").Value = -54989.5
# Row 94
$ws.Range("H94").Value = 59995
$ws.Range("J94").Value = 59995
$ws.Range("L94").Value = 59995
$ws.Range("N94").Value = -61797
# Row 97
$ws.Range("H97").Value = 206855.33
$ws.Range("J97").Value = 206855.33
$ws.Range("L97").Value = 206855.33
$ws.Range("N97").Value = -208837.33
# Row 132
$ws.Range("H132").Value = 16131467
$ws.Range("I132").Value = 20002002
$ws.Range("J132").Value = 4240
$ws.Range("K132").Value = 60006006
$ws.Range("L132").Value = 12720
$ws.Range("M132").Value = -60003476
$ws.Range("N132").Value = -17780
# Row 136
$ws.Range("H136").Value = 10786547
$ws.Range("I136").Value = 13931386
$ws.Range("K136").Value = 41794158
$ws.Range("M136").Value = -41791608
